$wb = $excel.ActiveWorkbook

# 1. Fix shared string / Tire_Type text: "710Rバフ50" -> "710R" in K2:K6 on Step3 sheets
foreach ($sheetName in @("Step3_DataPts_0.5","Step3_DataPts_0.7","Step3_DataPts_0.8","Step3_DataPts_0.9")) {
    $s = $wb.Worksheets.Item($sheetName)
    for ($r = 2; $r -le 6; $r++) {
        $s.Range("K$r").Value = "710R"
    }
}

# Step1_Data (69 cell updates)
$s = $wb.Worksheets.Item("Step1_Data")
$s.Range("D2").Value = 0.1846822601012179
$s.Range("E2").Value = 0.03612287545471059
$s.Range("F2").Value = 0.3421754535433115
$s.Range("G2").Value = 0.1505025059725814
$s.Range("M2").Value = 0.0002717033105636279
$s.Range("N2").Value = 0.01848680942741538
$s.Range("O2").Value = 0.1338826017316122
$s.Range("P2").Value = 0.05611776204902586
$s.Range("R2").Value = 0.004186786987433726
$s.Range("S2").Value = 0.01644636773052333
$s.Range("T2").Value = 0.00133399114960837
$s.Range("Z2").Value = 0.001622585697704025
$s.Range("AC2").Value = 0.005478047512542689
$s.Range("AD2").Value = 0.002480199335002122
$s.Range("AE2").Value = 0.01411555144945933
$s.Range("AH2").Value = 0.03209449854728772
$s.Range("D3").Value = 0.3670825490845082
$s.Range("F3").Value = 0.4588233992205751
$s.Range("G3").Value = 0.03668130700365096
$s.Range("H3").Value = 0.007784086107060921
$s.Range("I3").Value = 0.004877939741631982
$s.Range("O3").Value = 0.05628735673447342
$s.Range("P3").Value = 0.006825496494718258
$s.Range("T3").Value = 0.0220912063283072
$s.Range("Z3").Value = 0.001363599879531247
$s.Range("AB3").Value = 0.004379513555871962
$s.Range("AD3").Value = 0.002170179387268614
$s.Range("AE3").Value = 0.0134331378923263
$s.Range("AG3").Value = 0.01062593620549025
$s.Range("AH3").Value = 0.007574292364585204
$s.Range("D4").Value = 0.2502318013264258
$s.Range("F4").Value = 0.4560376009811742
$s.Range("G4").Value = 0.1328613007140707
$s.Range("I4").Value = 0.009337097654170715
$s.Range("O4").Value = 0.06133225956465835
$s.Range("P4").Value = 0.03512307784460952
$s.Range("S4").Value = 0.007455338003167065
$s.Range("T4").Value = 0.003301832398478458
$s.Range("Z4").Value = 0.001987832597119406
$s.Range("AD4").Value = 0.001695504342929497
$s.Range("AE4").Value = 0.02503426312104645
$s.Range("AH4").Value = 0.01560209145214949
$s.Range("D5").Value = 0.3678392964760315
$s.Range("E5").Value = 0.04683073161579145
$s.Range("F5").Value = 0.4408058585815912
$s.Range("G5").Value = 0.005518495593361934
$s.Range("H5").Value = 0.0143745827254378
$s.Range("I5").Value = 0.004984225349362353
$s.Range("O5").Value = 0.05628190251094038
$s.Range("R5").Value = 0.000428588743355322
$s.Range("T5").Value = 0.01849971645656474
$s.Range("Z5").Value = 0.003151335497256952
$s.Range("AB5").Value = 0.007284162925721574
$s.Range("AD5").Value = 0.009758427364976797
$s.Range("AE5").Value = 0.01528437457549021
$s.Range("AG5").Value = 0.008958301584117635
$s.Range("E6").Value = 0.2842842859730206
$s.Range("F6").Value = 0.09100225781906657
$s.Range("G6").Value = 0.4855514973602118
$s.Range("H6").Value = 0.01203494998421669
$s.Range("I6").Value = 0.02288420403985639
$s.Range("J6").Value = 0.01502498861760896
$s.Range("P6").Value = 0.03396743791293576
$s.Range("Q6").Value = 0.002973030301002553
$s.Range("U6").Value = 0.01079579158656727
$s.Range("AC6").Value = 0.01167705922835181
$s.Range("AE6").Value = 0.00699679190305257
$s.Range("AF6").Value = 0.01404168184336196
$s.Range("AH6").Value = 0.008766023430746882

# Step2_Sj (161 cell updates)
$s = $wb.Worksheets.Item("Step2_Sj")
$s.Range("D2").Value = 0.1846822601012179
$s.Range("E2").Value = 0.2208051355559285
$s.Range("F2").Value = 0.56298058909924
$s.Range("G2").Value = 0.7134830950718214
$s.Range("H2").Value = 0.7134830950718214
$s.Range("I2").Value = 0.7134830950718214
$s.Range("J2").Value = 0.7134830950718214
$s.Range("K2").Value = 0.7134830950718214
$s.Range("L2").Value = 0.7134830950718214
$s.Range("M2").Value = 0.713754798382385
$s.Range("N2").Value = 0.7322416078098004
$s.Range("O2").Value = 0.8661242095414127
$s.Range("P2").Value = 0.9222419715904385
$s.Range("Q2").Value = 0.9222419715904385
$s.Range("R2").Value = 0.9264287585778722
$s.Range("S2").Value = 0.9428751263083955
$s.Range("T2").Value = 0.944209117458004
$s.Range("U2").Value = 0.944209117458004
$s.Range("V2").Value = 0.944209117458004
$s.Range("W2").Value = 0.944209117458004
$s.Range("X2").Value = 0.944209117458004
$s.Range("Y2").Value = 0.944209117458004
$s.Range("Z2").Value = 0.945831703155708
$s.Range("AA2").Value = 0.945831703155708
$s.Range("AB2").Value = 0.945831703155708
$s.Range("AC2").Value = 0.9513097506682506
$s.Range("AD2").Value = 0.9537899500032527
$s.Range("AE2").Value = 0.967905501452712
$s.Range("AF2").Value = 0.967905501452712
$s.Range("AG2").Value = 0.967905501452712
$s.Range("AH2").Value = 0.9999999999999998
$s.Range("AI2").Value = 0.9999999999999998
$s.Range("AJ2").Value = 0.9999999999999998
$s.Range("D3").Value = 0.3670825490845082
$s.Range("E3").Value = 0.3670825490845082
$s.Range("F3").Value = 0.8259059483050833
$s.Range("G3").Value = 0.8625872553087343
$s.Range("H3").Value = 0.8703713414157952
$s.Range("I3").Value = 0.8752492811574272
$s.Range("J3").Value = 0.8752492811574272
$s.Range("K3").Value = 0.8752492811574272
$s.Range("L3").Value = 0.8752492811574272
$s.Range("M3").Value = 0.8752492811574272
$s.Range("N3").Value = 0.8752492811574272
$s.Range("O3").Value = 0.9315366378919006
$s.Range("P3").Value = 0.9383621343866189
$s.Range("Q3").Value = 0.9383621343866189
$s.Range("R3").Value = 0.9383621343866189
$s.Range("S3").Value = 0.9383621343866189
$s.Range("T3").Value = 0.9604533407149261
$s.Range("U3").Value = 0.9604533407149261
$s.Range("V3").Value = 0.9604533407149261
$s.Range("W3").Value = 0.9604533407149261
$s.Range("X3").Value = 0.9604533407149261
$s.Range("Y3").Value = 0.9604533407149261
$s.Range("Z3").Value = 0.9618169405944573
$s.Range("AA3").Value = 0.9618169405944573
$s.Range("AB3").Value = 0.9661964541503293
$s.Range("AC3").Value = 0.9661964541503293
$s.Range("AD3").Value = 0.9683666335375979
$s.Range("AE3").Value = 0.9817997714299241
$s.Range("AF3").Value = 0.9817997714299241
$s.Range("AG3").Value = 0.9924257076354144
$s.Range("AH3").Value = 0.9999999999999996
$s.Range("AI3").Value = 0.9999999999999996
$s.Range("AJ3").Value = 0.9999999999999996
$s.Range("D4").Value = 0.2502318013264258
$s.Range("E4").Value = 0.2502318013264258
$s.Range("F4").Value = 0.7062694023076
$s.Range("G4").Value = 0.8391307030216708
$s.Range("H4").Value = 0.8391307030216708
$s.Range("I4").Value = 0.8484678006758415
$s.Range("J4").Value = 0.8484678006758415
$s.Range("K4").Value = 0.8484678006758415
$s.Range("L4").Value = 0.8484678006758415
$s.Range("M4").Value = 0.8484678006758415
$s.Range("N4").Value = 0.8484678006758415
$s.Range("O4").Value = 0.9098000602404999
$s.Range("P4").Value = 0.9449231380851094
$s.Range("Q4").Value = 0.9449231380851094
$s.Range("R4").Value = 0.9449231380851094
$s.Range("S4").Value = 0.9523784760882765
$s.Range("T4").Value = 0.955680308486755
$s.Range("U4").Value = 0.955680308486755
$s.Range("V4").Value = 0.955680308486755
$s.Range("W4").Value = 0.955680308486755
$s.Range("X4").Value = 0.955680308486755
$s.Range("Y4").Value = 0.955680308486755
$s.Range("Z4").Value = 0.9576681410838744
$s.Range("AA4").Value = 0.9576681410838744
$s.Range("AB4").Value = 0.9576681410838744
$s.Range("AC4").Value = 0.9576681410838744
$s.Range("AD4").Value = 0.959363645426804
$s.Range("AE4").Value = 0.9843979085478504
$s.Range("AF4").Value = 0.9843979085478504
$s.Range("AG4").Value = 0.9843979085478504
$s.Range("D5").Value = 0.3678392964760315
$s.Range("E5").Value = 0.414670028091823
$s.Range("F5").Value = 0.8554758866734141
$s.Range("G5").Value = 0.860994382266776
$s.Range("H5").Value = 0.8753689649922138
$s.Range("I5").Value = 0.8803531903415762
$s.Range("J5").Value = 0.8803531903415762
$s.Range("K5").Value = 0.8803531903415762
$s.Range("L5").Value = 0.8803531903415762
$s.Range("M5").Value = 0.8803531903415762
$s.Range("N5").Value = 0.8803531903415762
$s.Range("O5").Value = 0.9366350928525166
$s.Range("P5").Value = 0.9366350928525166
$s.Range("Q5").Value = 0.9366350928525166
$s.Range("R5").Value = 0.9370636815958718
$s.Range("S5").Value = 0.9370636815958718
$s.Range("T5").Value = 0.9555633980524366
$s.Range("U5").Value = 0.9555633980524366
$s.Range("V5").Value = 0.9555633980524366
$s.Range("W5").Value = 0.9555633980524366
$s.Range("X5").Value = 0.9555633980524366
$s.Range("Y5").Value = 0.9555633980524366
$s.Range("Z5").Value = 0.9587147335496936
$s.Range("AA5").Value = 0.9587147335496936
$s.Range("AB5").Value = 0.9659988964754151
$s.Range("AC5").Value = 0.9659988964754151
$s.Range("AD5").Value = 0.9757573238403919
$s.Range("AE5").Value = 0.9910416984158822
$s.Range("AF5").Value = 0.9910416984158822
$s.Range("AG5").Value = 0.9999999999999998
$s.Range("AH5").Value = 0.9999999999999998
$s.Range("AI5").Value = 0.9999999999999998
$s.Range("AJ5").Value = 0.9999999999999998
$s.Range("E6").Value = 0.2842842859730206
$s.Range("F6").Value = 0.3752865437920872
$s.Range("G6").Value = 0.860838041152299
$s.Range("H6").Value = 0.8728729911365157
$s.Range("I6").Value = 0.8957571951763721
$s.Range("J6").Value = 0.910782183793981
$s.Range("K6").Value = 0.910782183793981
$s.Range("L6").Value = 0.910782183793981
$s.Range("M6").Value = 0.910782183793981
$s.Range("N6").Value = 0.910782183793981
$s.Range("O6").Value = 0.910782183793981
$s.Range("P6").Value = 0.9447496217069168
$s.Range("Q6").Value = 0.9477226520079194
$s.Range("R6").Value = 0.9477226520079194
$s.Range("S6").Value = 0.9477226520079194
$s.Range("T6").Value = 0.9477226520079194
$s.Range("U6").Value = 0.9585184435944867
$s.Range("V6").Value = 0.9585184435944867
$s.Range("W6").Value = 0.9585184435944867
$s.Range("X6").Value = 0.9585184435944867
$s.Range("Y6").Value = 0.9585184435944867
$s.Range("Z6").Value = 0.9585184435944867
$s.Range("AA6").Value = 0.9585184435944867
$s.Range("AB6").Value = 0.9585184435944867
$s.Range("AC6").Value = 0.9701955028228385
$s.Range("AD6").Value = 0.9701955028228385
$s.Range("AE6").Value = 0.9771922947258911
$s.Range("AF6").Value = 0.9912339765692531
$s.Range("AG6").Value = 0.9912339765692531
$s.Range("AH6").Value = 0.9999999999999999
$s.Range("AI6").Value = 0.9999999999999999
$s.Range("AJ6").Value = 0.9999999999999999

# Step3_DataPts_0.5 (5 cell updates)
$s = $wb.Worksheets.Item("Step3_DataPts_0.5")
$s.Range("F2").Value = 0.56298058909924
$s.Range("F3").Value = 0.8259059483050833
$s.Range("F4").Value = 0.7062694023076
$s.Range("F5").Value = 0.8554758866734141
$s.Range("F6").Value = 0.860838041152299

# Step3_DataPts_0.7 (5 cell updates)
$s = $wb.Worksheets.Item("Step3_DataPts_0.7")
$s.Range("F2").Value = 0.7134830950718214
$s.Range("F3").Value = 0.8259059483050833
$s.Range("F4").Value = 0.7062694023076
$s.Range("F5").Value = 0.8554758866734141
$s.Range("F6").Value = 0.860838041152299

# Step3_DataPts_0.8 (9 cell updates)
$s = $wb.Worksheets.Item("Step3_DataPts_0.8")
$s.Range("D2").Value = 14
$s.Range("F2").Value = 0.8661242095414127
$s.Range("G2").Value = 13
$s.Range("F3").Value = 0.8259059483050833
$s.Range("D4").Value = 6
$s.Range("F4").Value = 0.8391307030216708
$s.Range("G4").Value = 5
$s.Range("F5").Value = 0.8554758866734141
$s.Range("F6").Value = 0.860838041152299

# Step3_DataPts_0.9 (15 cell updates)
$s = $wb.Worksheets.Item("Step3_DataPts_0.9")
$s.Range("D2").Value = 15
$s.Range("F2").Value = 0.9222419715904385
$s.Range("G2").Value = 14
$s.Range("D3").Value = 14
$s.Range("F3").Value = 0.9315366378919006
$s.Range("G3").Value = 13
$s.Range("D4").Value = 14
$s.Range("F4").Value = 0.9098000602404999
$s.Range("G4").Value = 13
$s.Range("D5").Value = 14
$s.Range("F5").Value = 0.9366350928525166
$s.Range("G5").Value = 13
$s.Range("D6").Value = 9
$s.Range("F6").Value = 0.910782183793981
$s.Range("G6").Value = 7
